$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Data table (rows 16-22 -> 16-23): add a new worker (Carlos Javier Urbina
#    Chiquillo, 1 period) at the top of the detail rows, and add a new period
#    (2508) for the existing worker (Yoimer Oswaldo Hernandez Cuesta), who now
#    has periods 2502-2508 (ascending) instead of 2502-2507.
# ---------------------------------------------------------------------------

# Create the new last row (23) by copying the formatting (incl. the special
# bottom-border "last row" style) from the current last row (22).
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 16: new worker record (Carlos Javier Urbina Chiquillo)
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "1049930565"
$ws.Range("D16").Value2 = "CARLOS JAVIER URBINA CHIQUILLO"
$ws.Range("E16").Value2 = "2309"
$ws.Range("F16").Value2 = 15467
$ws.Range("G16").Value2 = 1300000

# Rows 17-23: Yoimer Oswaldo Hernandez Cuesta, periods 2502-2508 (ascending)
$periods = @("2502", "2503", "2504", "2505", "2506", "2507", "2508")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 17 + $i
    $ws.Range("B$r").Value2 = "CC"
    $ws.Range("C$r").Value2 = "1066576558"
    $ws.Range("D$r").Value2 = "YOIMER OSWALDO HERNANDEZ CUESTA"
    $ws.Range("E$r").Value2 = $periods[$i]
    $ws.Range("F$r").Value2 = 46400
    $ws.Range("G$r").Value2 = 1160000
}

# ---------------------------------------------------------------------------
# 2. Summary fields
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 340267   # Valor Mora (total)
$ws.Range("F13").Value2 = 8        # Cant. Periodos (7 -> 8)

# ---------------------------------------------------------------------------
# 3. Footer signature block: shift the two rows down by one so a blank row is
#    freed above them (rows 27-28 -> 28-29).
# ---------------------------------------------------------------------------
$ws.Rows("27:27").Insert()

# ---------------------------------------------------------------------------
# 4. Sheet dimension / used range now ends at row 29 instead of row 28.
# ---------------------------------------------------------------------------
Write-Output "edit complete"
